## Update navigation add media data
#
# - Fixes two stray "quote-prefix style" cells (B542, A557) back to the
#   plain text style used by the rest of the column.
# - Appends 8 new media rows (558-565) with Title / IMDB_ID / Type.
# - Moves the active selection down to the newly-added rows, matching the
#   navigation update described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize two legacy-formatted cells (style drifted during older edits) ---
$ws.Range("B542").Clear()
$ws.Range("B542").Value = "tt3506934"

$ws.Range("A557").Clear()
$ws.Range("A557").Value = "Okja"

# --- Append new media entries ---
$newRows = @(
    @("Wonder Woman", "tt0451279", "Movie"),
    @("Alien: Covenant", "tt2316204", "Movie"),
    @("My Bromance", "tt3522738", "Movie"),
    @("Ghost In The Shell", "tt1219827", "Movie"),
    @("Red Wine in the Dark Night", "tt4556730", "Movie"),
    @("The Mummy", "tt2345759", "Movie"),
    @("How to Win at Checkers (Every Time)", "tt4370256", "Movie"),
    @("Pink Moon", "tt3775450", "Movie")
)

$row = 558
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# --- Update navigation / selection to the tail of the new data ---
$ws.Range("A544").Select()
$ws.Range("C564:C565").Select()
